$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A currently holds a redundant index (1,2,4,13) duplicated from
# column F. Delete it so every other column shifts one place to the left,
# turning the old B:F data block into the new A:E block.
$ws.Range("A1:A5").Delete()
